# Auto-generated edit script: updates the "Price" (D) and composite
# "Volume(1h)" (E) text cells, plus the BKEXToken/KickToken row swap
# (B/C columns), to match the refreshed symbol-list snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.46"
$ws.Range("D3").Value = "'23.02"
$ws.Range("D4").Value = "'5.402"
$ws.Range("D6").Value = "'3.424"
$ws.Range("D7").Value = "'6.497"
$ws.Range("D8").Value = "'0.8129"
$ws.Range("D9").Value = "'0.9271"
$ws.Range("D10").Value = "'0.1439"
$ws.Range("D11").Value = "'0.07417"
$ws.Range("D12").Value = "'0.03318"
$ws.Range("D13").Value = "'0.03066"
$ws.Range("D14").Value = "'0.09350"
$ws.Range("D15").Value = "'3.848"
$ws.Range("D16").Value = "'0.001571"
$ws.Range("D17").Value = "'0.04711"
$ws.Range("D18").Value = "'0.0005909"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.005888"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("D22").Value = "'0.00007998"
$ws.Range("D23").Value = "'3.575"
$ws.Range("D27").Value = "'0.0002339"
$ws.Range("D40").Value = "'0.03953"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006449"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.003999"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1076"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.008565"
$ws.Range("D45").Value = "'0.00005171"
$ws.Range("D47").Value = "'0.6699"
$ws.Range("D48").Value = "'0.002273"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
